$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text, matching the
# original workbook convention where the Price column holds text values.

$ws.Range("D2").Value = "56.110.24"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").Value = "2.961.08"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("E4").Value = "  -1.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "500.68"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.93"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.430"
$ws.Range("E8").Value = "  +0.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.08"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  +0.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.367"
$ws.Range("E11").Value = "  +3.14%  "

$ws.Range("D12").Value = "3.467.92"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.125"
$ws.Range("E13").Value = "  -2.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.31"
$ws.Range("E14").Value = "  +1.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000162"
$ws.Range("E15").Value = "  +2.58%  "

$ws.Range("D16").Value = "55.585.78"
$ws.Range("E16").Value = "  -2.40%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.00"
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.951.57"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("E19").Value = "  +3.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.89"
$ws.Range("E20").Value = "  +1.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.60"
$ws.Range("E21").Value = "  +2.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.493"
$ws.Range("E23").Value = "  +2.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.89"
$ws.Range("E24").Value = "  +1.33%  "

$ws.Range("D25").Value = "3.105.15"
$ws.Range("E25").Value = "  +0.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  -0.83%  "

$ws.Range("D28").Value = "0.0₃0882"
$ws.Range("E28").Value = "  -0.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.50"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.04"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.78"
$ws.Range("E31").Value = "  +1.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("E32").Value = "  +1.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.25"
$ws.Range("E33").Value = "  +1.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.06"
$ws.Range("E34").Value = "  -2.53%  "

$ws.Range("E35").Value = "  -0.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.79"
$ws.Range("E36").Value = "  +1.41%  "

$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.31"
$ws.Range("E37").Value = "  +5.61%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.25"
$ws.Range("E38").Value = "  +0.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0661"
$ws.Range("E39").Value = "  -0.38%  "

$ws.Range("D40").Value = "2.963.62"
$ws.Range("E40").Value = "  -1.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.60"
$ws.Range("E41").Value = "  -2.42%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.79"
$ws.Range("E42").Value = "  +2.24%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.985"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.648"
$ws.Range("E44").Value = "  +1.98%  "

$ws.Range("D45").Value = "2.153.99"
$ws.Range("E45").Value = "  -2.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.36"
$ws.Range("E46").Value = "  -1.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.937"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.91"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.78"
$ws.Range("E49").Value = "  +3.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0236"
$ws.Range("E50").Value = "  +0.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0853"
$ws.Range("E51").Value = "  -2.29%  "
